$wb = $excel.ActiveWorkbook

# --- Reclassifications sheet: rename "MSCI UAE Index" to "MSCI UNITED ARAB EMIRATES Index" ---
$wsReclass = $wb.Worksheets.Item("Reclassifications")
$wsReclass.Cells.Item(10, 1).Value = "MSCI UNITED ARAB EMIRATES Index"

# --- Instructions sheet: add a new action row documenting the rename ---
$wsInstr = $wb.Worksheets.Item("Instructions")
$wsInstr.Cells.Item(4, 2).Value = '3) Rename "UAE" to "United Arab Emirates"'

# --- Update selections / active views ---
$wsReclass.Range("A10").Select()
$wsInstr.Activate()
$wsInstr.Range("B4").Select()
